# Loan RBI, Variable Instalments
# - Insert a new (empty) column before column N on the "Repayment schedule"
#   sheet, shifting the old N/O/P ("Late"/"Outstanding"/"Disbursement")
#   columns right to O/P/Q.
# - Make "Repayment schedule" the active sheet/tab, with K18 selected.

$wb = $excel.ActiveWorkbook

$wsSchedule = $wb.Worksheets.Item("Repayment schedule")

# Remember the width of column M (14 becomes the new column; M keeps its
# own width) so the freshly inserted column can be given the same width.
$mColumnWidth = $wsSchedule.Columns.Item(13).ColumnWidth

# Insert a blank column before column N (14th column); existing data in
# N:P shifts to O:Q automatically, carrying over values and styles.
$wsSchedule.Columns.Item(14).Insert()
$wsSchedule.Columns.Item(14).ColumnWidth = $mColumnWidth

# Activate the "Repayment schedule" sheet and move the selection to K18,
# matching the new tabSelected / active-cell state captured in the diff.
$wsSchedule.Activate()
$wsSchedule.Range("K18").Select()
